$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Insert the new "Const and Logical Const" Heading1 section, right after the
# paragraph ending "...then you do the copying of the data behind the
# object." (the end of the SharedByValue<T> discussion) and right before the
# "Synchronization (thread safety)" Heading1.
# ---------------------------------------------------------------------------

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like '*then you do the copying of the data behind the object.*') {
        $anchor = $cand
        break
    }
}

$anchorRange = $anchor.Range
$anchorIndex = $anchor.Range.Information(1)

# Create 6 fresh empty paragraphs immediately after the anchor paragraph.
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()

# Re-find the anchor paragraph's index now that new paragraphs exist, so we
# can reliably address the freshly inserted ones positionally.
$anchorPos = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like '*then you do the copying of the data behind the object.*') {
        $anchorPos = $i
        break
    }
}

$newHeading = $d.Paragraphs.Item($anchorPos + 1)
$newHeading.Style = 'Heading 1'
$newHeading.Range.Text = 'Const and Logical Const'

$newP1 = $d.Paragraphs.Item($anchorPos + 2)
$newP1.Range.Text = 'Generally Stroika uses the idea of logical const for its objects, and freely uses mutable for fields to enforce that notion.'

$newP2 = $d.Paragraphs.Item($anchorPos + 3)
$newP2.Range.Text = 'But there is one case where this is slightly vague, and at first glance, may appear not fully adhered to: Ptr objects.'

$newP3 = $d.Paragraphs.Item($anchorPos + 4)
$newP3.Range.Text = 'Ptr objects are really combinations of two kinds of things – smart pointers – and short-hand accessors for the underlying thing.'

$newP4 = $d.Paragraphs.Item($anchorPos + 5)
$newP4.Range.Text = 'Because of the C++ thread safety rules (always safe to access const methods from multiple threads at once so long as no writers, and the need for synchronization on writes) – and because these rules only apply literally and directly to the ‘envelope’ part – or the smart-pointer part of the object, we use the constness on Ptr objects to refer to the ptr itself, and not thing pointed to.'

$newP5 = $d.Paragraphs.Item($anchorPos + 6)
$newP5.Range.Text = 'We arguably COULD get rid of PTR objects and just use shared_ptr<T> or shared_ptr<const T> - but then we would lose the convenience of having simple interfaces for reps, and more complex, overloading etc interfaces for calling. '

# ---------------------------------------------------------------------------
# The trailing empty paragraph at the end of the document used to carry the
# (hidden) "_GoBack" last-edit-position bookmark. Now that real edits have
# been made earlier in the body, Word relocates that bookmark away from the
# final paragraph, leaving it a plain empty paragraph.
# ---------------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphBefore()

$lastIndex2 = $d.Paragraphs.Count
$oldLastPara = $d.Paragraphs.Item($lastIndex2)
$oldLastPara.Range.Delete()
